# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the newly generated gh-pages output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 167
$ws1.Range("F3").Value = 160
$ws1.Range("F5").Value = 4754
$ws1.Range("F8").Value = 523
$ws1.Range("F10").Value = 25
$ws1.Range("F12").Value = 1345
$ws1.Range("F13").Value = 2979
$ws1.Range("F14").Value = 387
$ws1.Range("F15").Value = 99
$ws1.Range("F16").Value = 86
$ws1.Range("F17").Value = 68
$ws1.Range("F18").Value = 2405
$ws1.Range("F19").Value = 107
$ws1.Range("F20").Value = 77
$ws1.Range("F22").Value = 168
$ws1.Range("F24").Value = 51
$ws1.Range("F25").Value = 240

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 167
$ws4.Range("F3").Value = 160
$ws4.Range("F6").Value = 4754
$ws4.Range("F9").Value = 523
$ws4.Range("F11").Value = 25
$ws4.Range("F13").Value = 1345
$ws4.Range("F14").Value = 2979
$ws4.Range("F15").Value = 387
$ws4.Range("F16").Value = 99
$ws4.Range("F17").Value = 86
$ws4.Range("F18").Value = 68
$ws4.Range("F19").Value = 2405
$ws4.Range("F20").Value = 107
$ws4.Range("F21").Value = 77
$ws4.Range("F23").Value = 168
$ws4.Range("F25").Value = 51
$ws4.Range("F26").Value = 240

$wb.Save()
